# Insert two new data rows (922 and 923) into the "Hortaliza, Macroferia Regional
# de Talca - Lechuga" sheet, pushing the existing rows 922:1006 down to 924:1008.
#
# Row 922 (new): Conconina(o) / Primera, dated 44578, "$/caja 10 unidades", Región del Maule
# Row 923 (new): Escarola / Primera, dated 44578, "$/caja 15 unidades", Región del Maule

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 922 (shifts 922:1006 -> 924:1008).
$ws.Rows.Item(922).Insert()
$ws.Rows.Item(922).Insert()

# Populate new row 922.
$ws.Cells.Item(922, 1).Value  = 5
$ws.Cells.Item(922, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(922, 3).Value  = "Maule"
$ws.Cells.Item(922, 4).Value  = 44578
$ws.Cells.Item(922, 5).Value  = 7
$ws.Cells.Item(922, 6).Value  = 100112033
$ws.Cells.Item(922, 7).Value  = "Lechuga"
$ws.Cells.Item(922, 8).Value  = "Conconina(o)"
$ws.Cells.Item(922, 9).Value  = "Primera"
$ws.Cells.Item(922, 10).Value = 500
$ws.Cells.Item(922, 11).Value = 4000
$ws.Cells.Item(922, 12).Value = 4000
$ws.Cells.Item(922, 13).Value = 4000
$ws.Cells.Item(922, 14).Value = "`$/caja 10 unidades"
$ws.Cells.Item(922, 15).Value = "Región del Maule"
$ws.Cells.Item(922, 16).Value = 400
$ws.Cells.Item(922, 17).Value = 10
$ws.Cells.Item(922, 18).Value = "Hortaliza"

# Populate new row 923.
$ws.Cells.Item(923, 1).Value  = 5
$ws.Cells.Item(923, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(923, 3).Value  = "Maule"
$ws.Cells.Item(923, 4).Value  = 44578
$ws.Cells.Item(923, 5).Value  = 7
$ws.Cells.Item(923, 6).Value  = 100112033
$ws.Cells.Item(923, 7).Value  = "Lechuga"
$ws.Cells.Item(923, 8).Value  = "Escarola"
$ws.Cells.Item(923, 9).Value  = "Primera"
$ws.Cells.Item(923, 10).Value = 450
$ws.Cells.Item(923, 11).Value = 6000
$ws.Cells.Item(923, 12).Value = 6000
$ws.Cells.Item(923, 13).Value = 6000
$ws.Cells.Item(923, 14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(923, 15).Value = "Región del Maule"
$ws.Cells.Item(923, 16).Value = 400
$ws.Cells.Item(923, 17).Value = 15
$ws.Cells.Item(923, 18).Value = "Hortaliza"
